{"js": "// Adicionado novas informa\u00e7oes nos dois arquivos\n//\n// 1) The title paragraph's single run is split into three runs, with a\n//    w:proofErr \"spellStart\"/\"spellEnd\" pair bracketing \"Git\" (what Word's\n//    proofer leaves behind for a word it doesn't recognise). The text itself\n//    is unchanged: \"Aula Git \u2013 Arquivo Local\".\n// 2) A new paragraph \"Thomas Novalski\" is added right after the title.\n// 3) A new, empty, trailing paragraph is added after that.\n\nconst body = context.document.body;\nbody.paragraphs.load(\"items\");\nawait context.sync();\n\n// --- 1) Re-run the title paragraph with proofErr markers around \"Git\" ---\n// Paragraph.insertOoxml requires a Flat OPC (\"pkg:package\") wrapper; build\n// one containing just the replacement <w:p> for this paragraph.\nconst titleFlatOpc = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t xml:space=\"preserve\">\\u201cAula </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>Git</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> \\u2013 Arquivo Local\\u201d</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nconst titleParagraph = body.paragraphs.items[0];\ntitleParagraph.insertOoxml(titleFlatOpc, \"Replace\");\nawait context.sync();\n\n// --- 2) Insert \"Thomas Novalski\" right after the (now re-run) title ---\nbody.paragraphs.load(\"items\");\nawait context.sync();\nconst titleAfterEdit = body.paragraphs.items[body.paragraphs.items.length - 1];\nconst namePara = titleAfterEdit.insertParagraph(\"Thomas Novalski\", \"After\");\nawait context.sync();\n\n// --- 3) Append one more, empty, trailing paragraph ---\n// insertParagraph always materializes a run, so insert a throwaway paragraph\n// and then replace it with a clean, run-less <w:p/> via insertOoxml.\nconst blankPara = namePara.insertParagraph(\"x\", \"After\");\nawait context.sync();\n\nconst blankFlatOpc = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body><w:p/></w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nblankPara.insertOoxml(blankFlatOpc, \"Replace\");\nawait context.sync();\n", "ps1": "# Adicionado novas informa\u00e7oes nos dois arquivos\n#\n# 1) The title paragraph's single run is split into three runs, with a\n#    w:proofErr \"spellStart\"/\"spellEnd\" pair bracketing \"Git\" (what Word's\n#    proofer leaves behind for a word it doesn't recognise). The text itself\n#    is unchanged: \"Aula Git \u2013 Arquivo Local\".\n# 2) A new paragraph \"Thomas Novalski\" is added right after the title.\n# 3) A new, empty, trailing paragraph is added after that.\n\n$d = $word.ActiveDocument\n\n# --- 1) Re-run the title paragraph with proofErr markers around \"Git\" ---\n# Range.InsertXML REPLACES the contents of the exact range it is called on,\n# so calling it on the title paragraph's own Range rewrites just that\n# paragraph's runs.\n$p1 = $d.Paragraphs.Item(1)\n$titleXml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:r><w:t xml:space=\"preserve\">\u201cAula </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Git</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> \u2013 Arquivo Local\u201d</w:t></w:r></w:p>'\n$p1.Range.InsertXML($titleXml)\n\n# --- 2) Insert \"Thomas Novalski\" right after the (now re-run) title ---\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n$namePara = $d.Paragraphs.Last\n$namePara.Range.Text = \"Thomas Novalski\"\n\n# --- 3) Append one more, empty, trailing paragraph ---\n# InsertParagraphAfter() always materializes an empty run (<w:r/>) in the new\n# paragraph; replace that paragraph's own Range content with a clean,\n# run-less <w:p/> via InsertXML (called on the paragraph's own Range, not a\n# collapsed end-of-document range, so only this paragraph is rewritten).\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n$blankPara = $d.Paragraphs.Last\n$blankPara.Range.InsertXML('<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"/>')\n"}
